$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.637506604194641
$ws.Range("B1").Value = 1.713898062705994
$ws.Range("C1").Value = 1.882205247879028
$ws.Range("D1").Value = 2.718746900558472
$ws.Range("E1").Value = 4.906844615936279
